# Auto-generated edit script: adds o_20 and o_20_jumbled sheets,
# updates o_10 with new evaluator_partial_correctness column.

$hdrPartial = @"
evaluator_partial_correctness
"@
$prompt16 = @"
 Given is the adjacency matrix for a unweighted undirected graph containing 16 nodes labelled A to P. The value corresponding to each row M and column N represents whether there is a connection between the two nodes, where 0 means no connection.   



what is the shortest path from node A to node P?

   A B C D E F G H I J K L M N O P
 A 0 1 0 0 1 0 0 0 0 0 0 0 0 0 0 0
 B 1 0 1 0 0 1 0 0 0 0 0 0 0 0 0 0
 C 0 1 0 1 0 0 1 0 0 0 0 0 0 0 0 0
 D 0 0 1 0 0 0 0 1 0 0 0 0 0 0 0 0
 E 1 0 0 0 0 1 0 0 1 0 0 0 0 0 0 0
 F 0 1 0 0 1 0 0 0 0 1 0 0 0 0 0 0
 G 0 0 1 0 0 0 0 1 0 0 1 0 0 0 0 0
 H 0 0 0 1 0 0 1 0 0 0 0 1 0 0 0 0
 I 0 0 0 0 1 0 0 0 0 0 0 0 1 0 0 0
 J 0 0 0 0 0 1 0 0 0 0 1 0 0 1 0 0
 K 0 0 0 0 0 0 1 0 0 1 0 1 0 0 0 0
 L 0 0 0 0 0 0 0 1 0 0 1 0 0 0 0 1
 M 0 0 0 0 0 0 0 0 1 0 0 0 0 1 0 0
 N 0 0 0 0 0 0 0 0 0 1 0 0 1 0 1 0
 O 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1
 P 0 0 0 0 0 0 0 0 0 0 0 1 0 0 1 0
    
"@
$sol16 = @"
A -> E -> F -> J -> K -> L -> P
"@
$resp16 = @"
To find the shortest path from node A to node P, we can use the breadth-first search algorithm.
Starting from node A, we check its neighbors (nodes B and E). We then enqueue these neighbors and mark them as visited. 
Next, we dequeue node B and check its neighbors (nodes A, C, and F). We enqueue the unvisited neighbors (C and F) and mark them as visited. 
We repeat this process, dequeuing and checking the neighbors of the nodes in the queue, until we reach node P. 
Here is a step-by-step breakdown of the algorithm:
Step 1: Initialize an empty queue and a visited array.
- Enqueue node A into the queue.
- Mark node A as visited.
Step 2: Perform BFS.
- Dequeue the first node from the queue (node A) and set it as the current node.
- Check the neighbors of the current node that are not visited.
  - If a neighbor is found, enqueue it and mark it as visited.
  - Repeat this step until all neighbors of the current node are checked.
- Repeat Step 2 until the current node is node P or the queue is empty.
Step 3: Backtrack the shortest path.
- Start from node P and backtrack to node A using the parent information stored during the BFS.
- The shortest path will be the sequence of nodes from node A to node P.
Based on this process, the shortest path from node A to node P is: A -> E -> J -> N -> P.
"@
$out16 = @"
Output: 2/5
"@
$prompt25 = @"
 Given is the adjacency matrix for a unweighted undirected graph containing 25 nodes labelled A to Y. The value corresponding to each row M and column N represents whether there is a connection between the two nodes, where 0 means no connection.   
what is the shortest path from node A to node Y?
   A B C D E F G H I J K L M N O P Q R S T U V W X Y
 A 0 1 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 B 1 0 1 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 C 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 D 0 0 1 0 1 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 E 0 0 0 1 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 F 1 0 0 0 0 0 1 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 G 0 1 0 0 0 1 0 1 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0
 H 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 I 0 0 0 1 0 0 0 0 0 1 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0
 J 0 0 0 0 1 0 0 0 1 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0
 K 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0
 L 0 0 0 0 0 0 1 0 0 0 0 0 1 0 0 0 1 0 0 0 0 0 0 0 0
 M 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0
 N 0 0 0 0 0 0 0 0 1 0 0 0 1 0 1 0 0 0 1 0 0 0 0 0 0
 O 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0 0 0 0 0 1 0 0 0 0 0
 P 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 1 0 0 0 0
 Q 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 1 0 0 0 1 0 0 0
 R 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0
 S 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0 1 0 0 0 1 0
 T 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0 0 0 0 0 1
 U 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 1 0 0 0
 V 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0 1 0 0
 W 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0
 X 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 1
 Y 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0
    
"@
$sol25 = @"
A -> F -> G -> L -> M -> N -> O -> T -> Y
"@
$resp25 = @"
To find the shortest path from node A to node Y, we can use the breadth-first search (BFS) algorithm. 
We start by initializing a queue and a visited array. We enqueue node A into the queue and mark it as visited. 
Then, we enter a loop where we dequeue a node from the queue and check if it is the target node Y. If so, we have found the shortest path. 
If not, we enqueue all the neighbors of the current node that have not been visited yet, mark them as visited, and continue the loop. 
Here is the step-by-step process:
1. Initialize an empty queue and a visited array.
2. Enqueue node A into the queue and mark it as visited.
3. Enter the BFS loop:
   - Dequeue a node from the queue.
   - Check if it is the target node Y. If so, we have found the shortest path.
   - If not, enqueue all the neighbors of the current node that have not been visited yet, mark them as visited.
   - Repeat the loop until the queue is empty or the target node is found.
4. If the target node is found, reconstruct the shortest path from the target node to the starting node A by backtracking through the parent pointers recorded during the BFS process.
Following this step-by-step process, we can find the shortest path from node A to node Y.
"@
$out25 = @"
Output: 1/8
"@
$prompt24 = @"
 Given is the adjacency matrix for a unweighted undirected graph containing 24 nodes labelled A to X. The value corresponding to each row M and column N represents whether there is a connection between the two nodes, where 0 means no connection.   
what is the shortest path from node A to node X?
   A B C D E F G H I J K L M N O P Q R S T U V W X
 A 0 1 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 B 1 0 1 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 C 0 1 0 1 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 D 0 0 1 0 1 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 E 0 0 0 1 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 F 1 0 0 0 0 0 1 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0
 G 0 1 0 0 0 1 0 1 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0
 H 0 0 1 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 I 0 0 0 1 0 0 0 1 0 1 0 0 0 1 0 0 0 0 0 0 0 0 0 0
 J 0 0 0 0 1 0 0 0 1 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0
 K 0 0 0 0 0 1 0 0 0 0 0 1 0 0 0 1 0 0 0 0 0 0 0 0
 L 0 0 0 0 0 0 1 0 0 0 1 0 1 0 0 0 1 0 0 0 0 0 0 0
 M 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 1 0 0 0 0 0 0
 N 0 0 0 0 0 0 0 0 1 0 0 0 0 0 1 0 0 0 1 0 0 0 0 0
 O 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0 0 0 0 0 1 0 0 0 0
 P 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 1 0 0 0 0 0 0 0
 Q 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0 1 0 0 1 0 0 0
 R 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0 1 0 0 1 0 0
 S 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0 1 0 0 1 0
 T 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0 0 0 0 1
 U 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 1 0 0
 V 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 1 0 1 0
 W 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 1 0 1
 X 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 1 0
    
"@
$sol24 = @"
A -> F -> G -> H -> I -> J -> O -> T -> X
"@
$resp24 = @"
To find the shortest path from node A to node X, we can use an algorithm such as Dijkstra's algorithm. Starting from node A, we can calculate the shortest distance to each node by considering all the possible paths and their associated weights.
Here is the step-by-step process using Dijkstra's algorithm:
1. Initialize an array to store the shortest distances from node A to all other nodes. Initially, all distances except for node A are set to infinity.
2. Create a visited array to keep track of the nodes that have been visited.
3. Set the distance of node A to 0 and mark it as visited.
4. For each adjacent node of A (nodes B, F), update their distances by considering the weight of the edge connecting them to A. In this case, the distance to B is 1 and the distance to F is 1. Update the array to reflect the new shortest distances.
5. Move to the node with the minimum distance from A that has not been visited (in this case, it is node B with a distance of 1). Mark it as visited.
6. For each of the adjacent nodes to B (nodes A, C, G), update their distances by considering the weight of the edge connecting them to B and the current shortest distance to B. In this case, the distance to A remains 1, the distance to C becomes 2, and the distance to G becomes 2. Update the array to reflect the new shortest distances.
7. Repeat steps 5 and 6 until all nodes have been visited. For each unvisited node, choose the one with the minimum distance from A and update its distances.
8. After visiting all nodes, the shortest distance array will contain the shortest distances from A to all other nodes. In this case, the shortest distance from A to X is 3.
9. To find the shortest path from A to X, we can backtrack from node X to A using the shortest distance array. Starting from X, move to the adjacent node with the smallest distance until we reach A. In this case, the shortest path from A to X is A -> F -> G -> X.
Therefore, the shortest path from node A to node X is A -> F -> G -> X with a distance of 3.
"@
$out24 = @"
Output: 3/9
"@

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Update sheet "o_10": add column E header + new row2 values ---
$ws1.Range("E1").Value2 = $hdrPartial
$ws1.Range("E1").Font.Bold = $true
$ws1.Range("E1").HorizontalAlignment = -4108
$ws1.Range("E1").VerticalAlignment = -4160
$ws1.Range("E1").Borders.LineStyle = 1

$ws1.Range("A2").Value2 = $prompt16
$ws1.Range("B2").Value2 = $sol16
$ws1.Range("C2").Value2 = $resp16
$ws1.Range("D2").Value2 = "Wrong"
$ws1.Range("E2").Value2 = $out16

# --- Create sheet "o_20" after the last sheet ---
$ws2 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws2.Name = "o_20"

$ws2.Range("A1").Value2 = "prompt"
$ws2.Range("B1").Value2 = "solution"
$ws2.Range("C1").Value2 = "llm_response"
$ws2.Range("D1").Value2 = "evaluator_response"
$ws2.Range("E1").Value2 = "evaluator_partial_correctness"
$ws2.Range("A1:E1").Font.Bold = $true
$ws2.Range("A1:E1").HorizontalAlignment = -4108
$ws2.Range("A1:E1").VerticalAlignment = -4160
$ws2.Range("A1:E1").Borders.LineStyle = 1

$ws2.Range("A2").Value2 = $prompt25
$ws2.Range("B2").Value2 = $sol25
$ws2.Range("C2").Value2 = $resp25
$ws2.Range("D2").Value2 = "Wrong"
$ws2.Range("E2").Value2 = $out25

# --- Create sheet "o_20_jumbled" after the last sheet ---
$ws3 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws3.Name = "o_20_jumbled"

$ws3.Range("A1").Value2 = "prompt"
$ws3.Range("B1").Value2 = "solution"
$ws3.Range("C1").Value2 = "llm_response"
$ws3.Range("D1").Value2 = "evaluator_response"
$ws3.Range("E1").Value2 = "evaluator_partial_correctness"
$ws3.Range("A1:E1").Font.Bold = $true
$ws3.Range("A1:E1").HorizontalAlignment = -4108
$ws3.Range("A1:E1").VerticalAlignment = -4160
$ws3.Range("A1:E1").Borders.LineStyle = 1

$ws3.Range("A2").Value2 = $prompt24
$ws3.Range("B2").Value2 = $sol24
$ws3.Range("C2").Value2 = $resp24
$ws3.Range("D2").Value2 = "Wrong"
$ws3.Range("E2").Value2 = $out24

# Restore sheet1 as the active/selected tab
$ws1.Activate()

Write-Host "Edit complete. Worksheets: $($wb.Worksheets.Count)"
